$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from O1 into the two new header columns P1 and Q1
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null

# Set the new header values (continuing the 0..15 sequence)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2-25: swap I<->K and M<->O values, and add new P/Q columns with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 10).Value = 2  # J (unchanged)
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 12).Value = 2  # L (unchanged)
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 14).Value = 2  # N (unchanged)
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q (new)
}
